$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A65").Value = "Both are fine."
$ws.Range("B65").Value = "どちらでもいいです。"
$ws.Range("A66").Value = "Same thing."
$ws.Range("B66").Value = "同じです。|おなじです。"
$ws.Range("A67").Value = "More or less the same."
$ws.Range("B67").Value = "だいたい同じです。|だいたいおなじです。"
$ws.Range("A68").Value = "A little different."
$ws.Range("B68").Value = "ちょっと違います。|ちょっとちがいます。"
$ws.Range("A69").Value = "Can't use it."
$ws.Range("B69").Value = "使えません。|つかえません。"
$ws.Range("A70").Value = "No good."
$ws.Range("B70").Value = "だめです。"
$ws.Range("A71").Value = "Raise your hand."
$ws.Range("B71").Value = "手をあげてください。|てをあげてください。"
$ws.Range("A72").Value = "Read it before coming to class."
$ws.Range("B72").Value = "読んできてください。|よんできてください。"
$ws.Range("A73").Value = "Hand in the homework."
$ws.Range("B73").Value = "宿題を出してください。|しゅくだいをだしてください。"
$ws.Range("A74").Value = "Open the book to page 10."
$ws.Range("B74").Value = "10ページを開いてください。|10ページをひらいてください。"
$ws.Range("A75").Value = "Close the textbook."
$ws.Range("B75").Value = "教科書を閉じてください。|きょうかしょをとじてください。"
$ws.Range("A76").Value = "Ask the person sitting next to you."
$ws.Range("B76").Value = "となりの人に聞いてください。|となりのひとにきいてください。"
$ws.Range("A77").Value = "Please stop."
$ws.Range("B77").Value = "やめてください。"
$ws.Range("A78").Value = "That's it for today."
$ws.Range("B78").Value = "今日はこれで終わります。|きょうはこれでおわります。"
$ws.Range("A79").Value = "deadline"
$ws.Range("B79").Value = "しめきり"
$ws.Range("A80").Value = "exercise"
$ws.Range("B80").Value = "練習|れんしゅう"
$ws.Range("A81").Value = "meaning"
$ws.Range("B81").Value = "意味|いみ"
$ws.Range("A82").Value = "pronunciation"
$ws.Range("B82").Value = "発音|はつおん"
$ws.Range("A83").Value = "grammar"
$ws.Range("B83").Value = "文法|ぶんぽう"
$ws.Range("A84").Value = "question"
$ws.Range("B84").Value = "質問|しつもん"
$ws.Range("A85").Value = "answer"
$ws.Range("B85").Value = "答|こたえ"
$ws.Range("A86").Value = "example"
$ws.Range("B86").Value = "例|れい"
$ws.Range("A87").Value = "( ) (parenthesis)"
$ws.Range("B87").Value = "かっこ"
$ws.Range("A88").Value = "〇 (correct)"
$ws.Range("B88").Value = "まる"
$ws.Range("A89").Value = "✕ (wrong)"
$ws.Range("B89").Value = "ばつ"
$ws.Range("A90").Value = "colloquial expression"
$ws.Range("B90").Value = "くだけた言い方|くだけたいいかた"
$ws.Range("A91").Value = "bookish expression"
$ws.Range("B91").Value = "かたい言い方|かたいいいかた"
$ws.Range("A92").Value = "polite expression"
$ws.Range("B92").Value = "ていねいな言い方|ていねいないいかた"
$ws.Range("A93").Value = "dialect"
$ws.Range("B93").Value = "方言|ほうげん"
$ws.Range("A94").Value = "common language"
$ws.Range("B94").Value = "共通語|きょうつうご"
$ws.Range("A95").Value = "for example"
$ws.Range("B95").Value = "たとえば"
$ws.Range("A96").Value = "anything else"
$ws.Range("B96").Value = "ほかに"
$ws.Range("A97").Value = "number"
$ws.Range("B97").Value = "～番|～ばん"
$ws.Range("A98").Value = "line number"
$ws.Range("B98").Value = "～行目|～ぎょうめ"
$ws.Range("A99").Value = "two people each"
$ws.Range("B99").Value = "二人ずつ|ふたりずつ"
